$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp column (Z) with new values from re-run of the notebook
$timestamps = @{
    2 = "2025-10-17T07:09:20.003943"
    3 = "2025-10-17T07:09:20.003943"
    4 = "2025-10-17T07:09:20.003943"
    5 = "2025-10-17T07:09:20.003943"
    6 = "2025-10-17T07:09:20.003943"
    7 = "2025-10-17T07:09:20.003943"
    8 = "2025-10-17T07:09:20.003943"
    9 = "2025-10-17T07:09:20.004945"
    10 = "2025-10-17T07:09:20.004945"
    11 = "2025-10-17T07:09:20.004945"
    12 = "2025-10-17T07:09:20.004945"
    13 = "2025-10-17T07:09:20.004945"
    14 = "2025-10-17T07:09:20.004945"
    15 = "2025-10-17T07:09:20.004945"
    16 = "2025-10-17T07:09:20.004945"
    17 = "2025-10-17T07:09:20.004945"
    18 = "2025-10-17T07:09:20.005943"
    19 = "2025-10-17T07:09:20.005943"
    20 = "2025-10-17T07:09:20.005943"
    21 = "2025-10-17T07:09:20.005943"
    22 = "2025-10-17T07:09:20.005943"
    23 = "2025-10-17T07:09:20.005943"
    24 = "2025-10-17T07:09:20.005943"
    25 = "2025-10-17T07:09:20.005943"
    26 = "2025-10-17T07:09:20.005943"
    27 = "2025-10-17T07:09:20.005943"
    28 = "2025-10-17T07:09:20.005943"
    29 = "2025-10-17T07:09:20.006943"
    30 = "2025-10-17T07:09:20.006943"
    31 = "2025-10-17T07:09:20.006943"
    32 = "2025-10-17T07:09:20.006943"
    33 = "2025-10-17T07:09:20.006943"
    34 = "2025-10-17T07:09:20.006943"
    35 = "2025-10-17T07:09:20.006943"
    36 = "2025-10-17T07:09:20.006943"
    37 = "2025-10-17T07:09:20.006943"
    38 = "2025-10-17T07:09:20.006943"
    39 = "2025-10-17T07:09:20.007943"
    40 = "2025-10-17T07:09:20.007943"
    41 = "2025-10-17T07:09:20.007943"
    42 = "2025-10-17T07:09:20.007943"
    43 = "2025-10-17T07:09:20.007943"
    44 = "2025-10-17T07:09:20.007943"
    45 = "2025-10-17T07:09:20.007943"
    46 = "2025-10-17T07:09:20.007943"
    47 = "2025-10-17T07:09:20.007943"
    48 = "2025-10-17T07:09:20.007943"
    49 = "2025-10-17T07:09:20.008943"
    50 = "2025-10-17T07:09:20.008943"
    51 = "2025-10-17T07:09:20.008943"
    52 = "2025-10-17T07:09:20.008943"
    53 = "2025-10-17T07:09:20.008943"
    54 = "2025-10-17T07:09:20.008943"
    55 = "2025-10-17T07:09:20.008943"
    56 = "2025-10-17T07:09:20.008943"
    57 = "2025-10-17T07:09:20.008943"
    58 = "2025-10-17T07:09:20.009943"
    59 = "2025-10-17T07:09:20.009943"
    60 = "2025-10-17T07:09:20.009943"
    61 = "2025-10-17T07:09:20.009943"
    62 = "2025-10-17T07:09:20.009943"
    63 = "2025-10-17T07:09:20.009943"
    64 = "2025-10-17T07:09:20.009943"
    65 = "2025-10-17T07:09:20.009943"
    66 = "2025-10-17T07:09:20.010941"
    67 = "2025-10-17T07:09:20.010941"
    68 = "2025-10-17T07:09:20.010941"
    69 = "2025-10-17T07:09:20.010941"
    70 = "2025-10-17T07:09:20.010941"
    71 = "2025-10-17T07:09:20.010941"
    72 = "2025-10-17T07:09:20.010941"
    73 = "2025-10-17T07:09:20.010941"
    74 = "2025-10-17T07:09:20.011942"
    75 = "2025-10-17T07:09:20.011942"
    76 = "2025-10-17T07:09:20.011942"
    77 = "2025-10-17T07:09:20.012943"
    78 = "2025-10-17T07:09:20.012943"
    79 = "2025-10-17T07:09:20.012943"
    80 = "2025-10-17T07:09:20.012943"
    81 = "2025-10-17T07:09:20.012943"
    82 = "2025-10-17T07:09:20.012943"
    83 = "2025-10-17T07:09:20.012943"
    84 = "2025-10-17T07:09:20.013942"
    85 = "2025-10-17T07:09:20.013942"
    86 = "2025-10-17T07:09:20.013942"
    87 = "2025-10-17T07:09:20.013942"
    88 = "2025-10-17T07:09:20.013942"
    89 = "2025-10-17T07:09:20.013942"
    90 = "2025-10-17T07:09:20.013942"
    91 = "2025-10-17T07:09:20.013942"
    92 = "2025-10-17T07:09:20.013942"
    93 = "2025-10-17T07:09:20.014942"
    94 = "2025-10-17T07:09:20.014942"
    95 = "2025-10-17T07:09:20.014942"
    96 = "2025-10-17T07:09:20.014942"
    97 = "2025-10-17T07:09:20.014942"
    98 = "2025-10-17T07:09:20.014942"
    99 = "2025-10-17T07:09:20.014942"
    100 = "2025-10-17T07:09:20.014942"
    101 = "2025-10-17T07:09:20.014942"
    102 = "2025-10-17T07:09:20.015940"
    103 = "2025-10-17T07:09:20.015940"
    104 = "2025-10-17T07:09:20.015940"
    105 = "2025-10-17T07:09:20.015940"
    106 = "2025-10-17T07:09:20.015940"
    107 = "2025-10-17T07:09:20.015940"
    108 = "2025-10-17T07:09:20.015940"
    109 = "2025-10-17T07:09:20.015940"
    110 = "2025-10-17T07:09:20.015940"
    111 = "2025-10-17T07:09:20.015940"
    112 = "2025-10-17T07:09:20.016941"
    113 = "2025-10-17T07:09:20.016941"
    114 = "2025-10-17T07:09:20.016941"
    115 = "2025-10-17T07:09:20.016941"
    116 = "2025-10-17T07:09:20.016941"
    117 = "2025-10-17T07:09:20.016941"
    118 = "2025-10-17T07:09:20.016941"
    119 = "2025-10-17T07:09:20.016941"
    120 = "2025-10-17T07:09:20.016941"
    121 = "2025-10-17T07:09:20.017940"
    122 = "2025-10-17T07:09:20.017940"
    123 = "2025-10-17T07:09:20.017940"
    124 = "2025-10-17T07:09:20.017940"
    125 = "2025-10-17T07:09:20.017940"
    126 = "2025-10-17T07:09:20.017940"
    127 = "2025-10-17T07:09:20.017940"
    128 = "2025-10-17T07:09:20.018940"
    129 = "2025-10-17T07:09:20.018940"
    130 = "2025-10-17T07:09:20.018940"
    131 = "2025-10-17T07:09:20.018940"
    132 = "2025-10-17T07:09:20.018940"
    133 = "2025-10-17T07:09:20.018940"
    134 = "2025-10-17T07:09:20.018940"
    135 = "2025-10-17T07:09:20.018940"
    136 = "2025-10-17T07:09:20.018940"
    137 = "2025-10-17T07:09:20.019940"
    138 = "2025-10-17T07:09:20.019975"
    139 = "2025-10-17T07:09:20.019975"
    140 = "2025-10-17T07:09:20.019975"
    141 = "2025-10-17T07:09:20.019975"
    142 = "2025-10-17T07:09:20.019975"
    143 = "2025-10-17T07:09:20.019975"
    144 = "2025-10-17T07:09:20.019975"
    145 = "2025-10-17T07:09:20.019975"
    146 = "2025-10-17T07:09:20.019975"
    147 = "2025-10-17T07:09:20.019975"
    148 = "2025-10-17T07:09:20.019975"
    149 = "2025-10-17T07:09:20.019975"
    150 = "2025-10-17T07:09:20.019975"
    151 = "2025-10-17T07:09:20.019975"
    152 = "2025-10-17T07:09:20.019975"
    153 = "2025-10-17T07:09:20.019975"
    154 = "2025-10-17T07:09:20.019975"
    155 = "2025-10-17T07:09:20.019975"
    156 = "2025-10-17T07:09:20.019975"
    157 = "2025-10-17T07:09:20.019975"
    158 = "2025-10-17T07:09:20.019975"
    159 = "2025-10-17T07:09:20.019975"
    160 = "2025-10-17T07:09:20.019975"
    161 = "2025-10-17T07:09:20.019975"
    162 = "2025-10-17T07:09:20.019975"
    163 = "2025-10-17T07:09:20.019975"
    164 = "2025-10-17T07:09:20.019975"
    165 = "2025-10-17T07:09:20.019975"
    166 = "2025-10-17T07:09:20.019975"
    167 = "2025-10-17T07:09:20.019975"
    168 = "2025-10-17T07:09:20.019975"
    169 = "2025-10-17T07:09:20.019975"
    170 = "2025-10-17T07:09:20.019975"
    171 = "2025-10-17T07:09:20.019975"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
